# first of many passes through in cleaning database
# Column I (roboticS1Prep) was text "No" for every data row; convert it to a
# real boolean FALSE value (with a TRUE/FALSE custom display format), which
# also drops the now-unused "No" shared string and renumbers "random".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("I2:I41")
$rng.Value = $false
$rng.NumberFormat = """TRUE"";""TRUE"";""FALSE"""

# Move the active selection to the column that was just edited, and scroll
# the view down a bit, matching where the author was working.
$excel.Goto($ws.Range("A7"), $true)
$ws.Range("I2:I41").Select()
